$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.636.10'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '2.531.72'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Formula = '="593.76"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Formula = '="177.52"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +2.10%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('D9').Value = '2.532.21'
$ws.Range('E9').Value = '  +2.66%  '
$ws.Range('D10').Formula = '="0.146"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  +6.40%  '
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').Value = '68.516.32'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '2.529.26'
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').Formula = '="11.11"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  +2.08%  '
$ws.Range('D20').Formula = '="7.54"'
$ws.Range('D20').Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').Formula = '="353.35"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('E22').Value = '  +4.99%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Formula = '="70.95"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').Formula = '="4.22"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('D26').Formula = '="1.71"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -4.43%  '
$ws.Range('D27').Formula = '="9.03"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('D28').Value = '2.656.61'
$ws.Range('E28').Value = '  +2.33%  '
$ws.Range('D29').Formula = '="1.00"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('D30').Formula = '="513.05"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  +3.15%  '
$ws.Range('D31').Value = '0.0₃0899'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('E33').Value = '  +2.77%  '
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Formula = '="164.13"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Formula = '="18.43"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  +1.82%  '
$ws.Range('D40').Formula = '="1.32"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  +4.69%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').Formula = '="0.327"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').Formula = '="2.43"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  +2.04%  '
$ws.Range('D46').Formula = '="152.70"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +7.34%  '
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Formula = '="0.522"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0261'
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('E51').Value = '  +0.64%  '

$excel.CutCopyMode = 0

